{"js": "// The document repeats a 2018 Perseus observation-dates sentence in four\n// paragraphs; each paragraph currently holds several differently-formatted\n// runs (\"2018 \" bold-italic + plain-text runs). The edit collapses every\n// such paragraph down to a single, unformatted run with the new Gemini\n// wording (\"\u0397\u03bc\u03b5\u03c1\u03bf\u03bc\u03b7\u03bd\u03af\u03b5\u03c2 \u03c0\u03b1\u03c1\u03b1\u03c4\u03ae\u03c1\u03b7\u03c3\u03b7\u03c2 \u03b3\u03b9\u03b1 \u03c4\u03bf\u03bd \u03b1\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc \u03c4\u03bf\u03c5 Gemini: 14-23\n// \u03a6\u03b5\u03b2\u03c1\u03bf\u03c5\u03b1\u03c1\u03af\u03bf\u03c5, 14-24 \u039c\u03b1\u03c1\u03c4\u03af\u03bf\u03c5\").\n\nconst newText =\n  \"\u0397\u03bc\u03b5\u03c1\u03bf\u03bc\u03b7\u03bd\u03af\u03b5\u03c2 \u03c0\u03b1\u03c1\u03b1\u03c4\u03ae\u03c1\u03b7\u03c3\u03b7\u03c2 \u03b3\u03b9\u03b1 \u03c4\u03bf\u03bd \u03b1\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc \u03c4\u03bf\u03c5 Gemini: 14-23 \u03a6\u03b5\u03b2\u03c1\u03bf\u03c5\u03b1\u03c1\u03af\u03bf\u03c5, 14-24 \u039c\u03b1\u03c1\u03c4\u03af\u03bf\u03c5\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Identify every paragraph whose text begins with the old \"2018\n// \u0397\u03bc\u03b5\u03c1\u03bf\u03bc\u03b7\u03bd\u03af\u03b5\u03c2...\" sentence (there are 4 of them in this document; a 5th,\n// unrelated paragraph also mentions \"\u03a0\u03b5\u03c1\u03c3\u03b5\u03cd\u03c2\" in running prose and must be\n// left untouched).\nconst targets = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"2018 \u0397\u03bc\u03b5\u03c1\u03bf\u03bc\u03b7\u03bd\u03af\u03b5\u03c2\") === 0) {\n    targets.push(paragraphs.items[i]);\n  }\n}\n\nfor (const paragraph of targets) {\n  // Remove every existing (variously formatted) run in the paragraph...\n  paragraph.clear();\n  // ...and replace it with a single plain run holding the new text.\n  paragraph.insertText(newText, \"Start\");\n}\n\nawait context.sync();\n", "ps1": "# The document repeats a 2018 Perseus observation-dates sentence in four\n# paragraphs; each paragraph currently holds several differently-formatted\n# runs (\"2018 \" bold-italic + plain-text runs). The edit collapses every\n# such paragraph down to a single, unformatted run with the new Gemini\n# wording (\"\u0397\u03bc\u03b5\u03c1\u03bf\u03bc\u03b7\u03bd\u03af\u03b5\u03c2 \u03c0\u03b1\u03c1\u03b1\u03c4\u03ae\u03c1\u03b7\u03c3\u03b7\u03c2 \u03b3\u03b9\u03b1 \u03c4\u03bf\u03bd \u03b1\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc \u03c4\u03bf\u03c5 Gemini: 14-23\n# \u03a6\u03b5\u03b2\u03c1\u03bf\u03c5\u03b1\u03c1\u03af\u03bf\u03c5, 14-24 \u039c\u03b1\u03c1\u03c4\u03af\u03bf\u03c5\").\n\n$d = $word.ActiveDocument\n\n$needle = \"2018 \u0397\u03bc\u03b5\u03c1\u03bf\u03bc\u03b7\u03bd\u03af\u03b5\u03c2\"\n$newText = \"\u0397\u03bc\u03b5\u03c1\u03bf\u03bc\u03b7\u03bd\u03af\u03b5\u03c2 \u03c0\u03b1\u03c1\u03b1\u03c4\u03ae\u03c1\u03b7\u03c3\u03b7\u03c2 \u03b3\u03b9\u03b1 \u03c4\u03bf\u03bd \u03b1\u03c3\u03c4\u03b5\u03c1\u03b9\u03c3\u03bc\u03cc \u03c4\u03bf\u03c5 Gemini: 14-23 \u03a6\u03b5\u03b2\u03c1\u03bf\u03c5\u03b1\u03c1\u03af\u03bf\u03c5, 14-24 \u039c\u03b1\u03c1\u03c4\u03af\u03bf\u03c5\"\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n\n    # Only the four paragraphs that still carry the old \"2018 \u0397\u03bc\u03b5\u03c1\u03bf\u03bc\u03b7\u03bd\u03af\u03b5\u03c2...\"\n    # sentence are rewritten; a fifth, unrelated paragraph also mentions\n    # \"\u03a0\u03b5\u03c1\u03c3\u03b5\u03cd\u03c2\" in running prose and must stay untouched.\n    if ($t.Length -ge $needle.Length -and $t.Substring(0, $needle.Length) -eq $needle) {\n        $r = $p.Range\n        # Exclude the trailing paragraph mark from the range so only the\n        # sentence's text (and its runs) is affected.\n        $r.MoveEnd(1, -1) | Out-Null\n        # Delete every existing (variously formatted) run in the paragraph...\n        $r.Delete()\n        # ...then insert a single new, unformatted run with the new text.\n        $r.InsertAfter($newText)\n    }\n}\n"}
